# Insert a new weekly price record before the existing row 534.
# This pushes the current rows 534:588 down to 535:589 and the
# worksheet dimension grows from A1:R588 to A1:R589.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(534).Insert()

$ws.Cells.Item(534, 1).Value  = 8
$ws.Cells.Item(534, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(534, 3).Value  = "Coquimbo"
$ws.Cells.Item(534, 4).Value  = 45194
$ws.Cells.Item(534, 5).Value  = 4
$ws.Cells.Item(534, 6).Value  = 100112032
$ws.Cells.Item(534, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(534, 8).Value  = "Sin especificar"
$ws.Cells.Item(534, 9).Value  = "Primera"
$ws.Cells.Item(534, 10).Value = 400
$ws.Cells.Item(534, 11).Value = 12000
$ws.Cells.Item(534, 12).Value = 13000
$ws.Cells.Item(534, 13).Value = 12500
$ws.Cells.Item(534, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(534, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(534, 16).Value = 250
$ws.Cells.Item(534, 17).Value = 50
$ws.Cells.Item(534, 18).Value = "Hortaliza"
